# pruebaV2.xlsx update — "Add files via upload"
#
# The author re-worded three step-description cells (switching the verb
# from first-person-plural "-amos" phrasing to the imperative form used
# everywhere else in the table) and left the selection sitting on B17
# (the row they'd just edited) instead of the previous G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (paso 17): "Tomamos el punto $G$" -> "Tomar el punto $G$"
$ws.Range("B18").Value = "Tomar el punto `$G`$"

# Row 8 (paso 7): "Prolongamos el segmento ..." -> "Prolongar el segmento ..."
$ws.Range("B8").Value = "Prolongar el segmento `$\overline{BC}`$ a un punto `$D`$"

# Row 17 (paso 16): "Quitamos del segmento ..." -> "Quitar del segmento ..."
$ws.Range("B17").Value = "Quitar del segmento `$\overline{EF}`$ el segmento `$\overline{BE}`$ obteniendo el segmento `$\overline{GF}`$"

# Selection moved from G7 to B17
$ws.Range("B17").Select()
